# Ruchi Pareek - "Changed deck"
# Slide 1, TextBox 1: split the existing bullet's wording across two runs
# and add a second auto-numbered bullet describing the edit session. The
# shape uses <a:spAutoFit/>, so re-entering the text through the normal
# TextRange/InsertAfter COM calls lets PowerPoint recompute the textbox
# height for the extra line, just like the author's edit did.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# First bullet: "Skill Gap introduce a new perspective" typed back in as
# two runs ("...a new " / "perspective") - this also drops the old
# trailing space after "perspective".
$tr.Text = "Skill Gap introduce a new "
$run1b = $tr.InsertAfter("perspective")

# New second bullet (same auto-numbered list - a new paragraph inherits
# the list's pPr automatically): "Editing script.pptx from Ruchi's system "
$apostrophe = [char]0x2019
$run2a = $run1b.InsertAfter("`rEditing script.pptx from ")
$run2b = $run2a.InsertAfter("Ruchi" + $apostrophe + "s")
$run2c = $run2b.InsertAfter(" system")
$run2d = $run2c.InsertAfter(" ")
